$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B8").Value = "2025-10-02T18:31:12+01:00"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "true"
